# Added a method to get permission list by ID number.
# Also disgusting support for remapping permissions.
# Permission remap only works with add_person, if setting or
# removing permission only the ones called will have an effect.

$wb = $excel.ActiveWorkbook

# "clearance" sheet gets a new "Top-secret" permission column (C) and
# becomes the active/selected sheet (mirrors the author switching to it
# while working on the remap feature).
$clearance = $wb.Worksheets.Item("clearance")
$people = $wb.Worksheets.Item("people")

$clearance.Range("C1").Value = "Top-secret"
$clearance.Range("C2").Value = 123

# Select C1 on the clearance sheet and make it the active tab.
$clearance.Activate()
$clearance.Range("C1").Select()
